# Applies the commit's text edits:
#  - Slide 17 title: fix typo "Red Ablation" -> "Redo Ablation" (splits into 3 runs,
#    matching the selection-retype pattern "of Red " -> "of Redo ")
#  - Slide 21 "Summary" bullets:
#      * "Base" row: add an extra tab before "+2"
#      * "PVI" row: add an extra tab before "+5"
#      * "Age >= 70 years" row: add an extra tab before "- 1" (splits into 3 runs)

$p = $ppt.ActivePresentation

# --- Slide 17: title typo fix -------------------------------------------------
$s17 = $p.Slides.Item(17)
$titleShape = $s17.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleFull = $titleRange.Text
$idx = $titleFull.IndexOf("of Red ")
if ($idx -ge 0) {
    $sub = $titleRange.Characters($idx + 1, 7)
    $sub.Text = "of Redo "
}

# --- Slide 21: scoring list tab fixes ----------------------------------------
$s21 = $p.Slides.Item(21)
$bodyShape = $s21.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange

# "Base" row: "Base\t\t\t+2" -> "Base\t\t\t\t+2"
$full = $bodyRange.Text
$idx = $full.IndexOf("Base`t`t`t+2")
if ($idx -ge 0) {
    $sub = $bodyRange.Characters($idx + 1, 9)
    $sub.Text = "Base`t`t`t`t+2"
}

# "PVI" row: "PVI\t\t\t+5" -> "PVI\t\t\t\t+5"
$full = $bodyRange.Text
$idx = $full.IndexOf("PVI`t`t`t+5")
if ($idx -ge 0) {
    $sub = $bodyRange.Characters($idx + 1, 8)
    $sub.Text = "PVI`t`t`t`t+5"
}

# "Age >= 70 years" row: "Age >= 70 years\t- 1" -> "Age >= 70 years\t\t- 1"
# (retype just the "\t- " portion so "Age >= 70 years" and the trailing "1"
#  stay in their own runs, matching the select+retype pattern of the edit)
$full = $bodyRange.Text
$idx = $full.IndexOf("Age >= 70 years")
if ($idx -ge 0) {
    $start = $idx + 1 + "Age >= 70 years".Length
    $sub = $bodyRange.Characters($start, 3)
    $sub.Text = "`t`t- "
}
